# ------------------------------------------------------------------
# B6-PowerPoint.pptx edit:
#   1) Three tables (slides 14, 15, 16) get their table style switched
#      from the custom "Table_0" style to the built-in table style
#      {9786098F-5CC9-4F8F-B038-4F61F09C066B}.
#   2) The slide theme (ppt/theme/theme1.xml, used by the slide master)
#      and the notes theme (ppt/theme/theme2.xml, used by the notes
#      master) swap their 12-colour theme colour schemes: the slides
#      now use the stock "Office" palette while the notes pages use
#      the palette that used to belong to the slides ("Red Violet").
# ------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Retarget the three tables' style ---------------------------
$newTableStyle = "{9786098F-5CC9-4F8F-B038-4F61F09C066B}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2) Swap the theme colour schemes -------------------------------
# Office theme colours (what the slide theme should become)
$officeColors = @{
    1  = 0         # dk1      000000
    2  = 16777215  # lt1      FFFFFF
    3  = 6968388   # dk2      44546A
    4  = 15132391  # lt2      E7E6E6
    5  = 13998939  # accent1  5B9BD5
    6  = 3243501    # accent2  ED7D31
    7  = 10855845   # accent3  A5A5A5
    8  = 49407      # accent4  FFC000
    9  = 12874308   # accent5  4472C4
    10 = 4697456    # accent6  70AD47
    11 = 12673797   # hlink    0563C1
    12 = 7491477    # folHlink 954F72
}

# Red Violet / Integral theme colours (what the notes theme should become)
$redVioletColors = @{
    1  = 0         # dk1      000000
    2  = 16777215  # lt1      FFFFFF
    3  = 5326149   # dk2      454551
    4  = 14473688  # lt2      D8D9DC
    5  = 9514467   # accent1  E32D91
    6  = 13381832  # accent2  C830CC
    7  = 14460494  # accent3  4EA6DC
    8  = 15168839  # accent4  4775E7
    9  = 14774665  # accent5  8971E1
    10 = 7555029   # accent6  D54773
    11 = 2465643   # hlink    6B9F25
    12 = 9211020   # folHlink 8C8C8C
}

$slide1 = $p.Slides.Item(1)

# theme1.xml backs the slide master -> drives the slides themselves
$slideThemeColors = $slide1.ThemeColorScheme
for ($k = 1; $k -le $slideThemeColors.Count; $k++) {
    $slideThemeColors.Item($k).RGB = $officeColors[$k]
}

# theme2.xml backs the notes master -> drives the notes pages
$notesThemeColors = $slide1.NotesPage.ThemeColorScheme
for ($k = 1; $k -le $notesThemeColors.Count; $k++) {
    $notesThemeColors.Item($k).RGB = $redVioletColors[$k]
}
